# Updates cryptos list values (price + 1h volume %) and reorders two coin rows
# per the scraper run on Tue Dec 12 03:47:33 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23/24: ImmutableX and BitcoinCash swap ranking positions
# Row 29/30: Toncoin and Monero swap ranking positions
# D-column price text is prefixed with a literal apostrophe so Excel keeps it as
# text (matching the source sheet, which stores prices/volumes as inline strings)
# instead of auto-converting look-alike numbers like "252.59" to a numeric value.

$ws.Range("D2").Value = '''41.859.98'
$ws.Range("E2").Value = '  -1.58%  '

$ws.Range("D3").Value = '''2.233.48'
$ws.Range("E3").Value = '  -0.98%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").Value = '''252.59'
$ws.Range("E5").Value = '  +8.17%  '

$ws.Range("D6").Value = '''0.625'
$ws.Range("E6").Value = '  -2.21%  '

$ws.Range("D7").Value = '''71.59'
$ws.Range("E7").Value = '  -0.75%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").Value = '''0.566'
$ws.Range("E9").Value = '  -0.51%  '

$ws.Range("D10").Value = '''42.09'
$ws.Range("E10").Value = '  +13.86%  '

$ws.Range("D11").Value = '''0.0971'
$ws.Range("E11").Value = '  -2.59%  '

$ws.Range("D12").Value = '''58.48'
$ws.Range("E12").Value = '  -0.74%  '

$ws.Range("D13").Value = '''0.105'
$ws.Range("E13").Value = '  -0.14%  '

$ws.Range("E14").Value = '  +0.93%  '

$ws.Range("D15").Value = '''2.559.96'
$ws.Range("E15").Value = '  -1.23%  '

$ws.Range("D16").Value = '''15.06'
$ws.Range("E16").Value = '  -0.71%  '

$ws.Range("D17").Value = '''0.860'
$ws.Range("E17").Value = '  -2.26%  '

$ws.Range("D18").Value = '''2.238.81'
$ws.Range("E18").Value = '  -0.49%  '

$ws.Range("D19").Value = '''41.756.34'
$ws.Range("E19").Value = '  -1.61%  '

$ws.Range("E20").Value = '  -2.34%  '

$ws.Range("D21").Value = '''73.17'
$ws.Range("E21").Value = '  -1.16%  '

$ws.Range("E22").Value = '  -1.04%  '

$ws.Range("B23").Value = 'ImmutableX'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D23").Value = '''2.25'
$ws.Range("E23").Value = '  +15.04%  '

$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").Value = '''235.27'
$ws.Range("E24").Value = '  -1.55%  '

$ws.Range("E25").Value = '  -0.06%  '

$ws.Range("D26").Value = '''3.76'
$ws.Range("E26").Value = '  +1.78%  '

$ws.Range("D27").Value = '''2.50'
$ws.Range("E27").Value = '  +5.27%  '

$ws.Range("D28").Value = '''10.25'
$ws.Range("E28").Value = '  +0.83%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '''2.20'
$ws.Range("E29").Value = '  +1.18%  '

$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = '''170.22'
$ws.Range("E30").Value = '  +1.03%  '

$ws.Range("D31").Value = '''20.79'
$ws.Range("E31").Value = '  -0.44%  '

$ws.Range("E32").Value = '  -0.15%  '

$ws.Range("E33").Value = '  -2.17%  '

$ws.Range("D34").Value = '''5.46'
$ws.Range("E34").Value = '  +0.83%  '

$ws.Range("D35").Value = '''0.0722'
$ws.Range("E35").Value = '  -1.26%  '

$ws.Range("D36").Value = '''26.81'
$ws.Range("E36").Value = '  +18.05%  '

$ws.Range("D37").Value = '''4.67'
$ws.Range("E37").Value = '  -3.49%  '

$ws.Range("E38").Value = '  +12.19%  '

$ws.Range("D39").Value = '''0.0283'
$ws.Range("E39").Value = '  +5.61%  '

$ws.Range("E40").Value = '  +1.13%  '

$ws.Range("D41").Value = '''69.73'

$ws.Range("D42").Value = '''6.04'
$ws.Range("E42").Value = '  -2.08%  '

$ws.Range("E43").Value = '  +11.57%  '

$ws.Range("D44").Value = '''5.09'
$ws.Range("E44").Value = '  -1.96%  '

$ws.Range("D45").Value = '''11.68'
$ws.Range("E45").Value = '  +13.06%  '

$ws.Range("D46").Value = '''8.87'
$ws.Range("E46").Value = '  -1.11%  '

$ws.Range("D47").Value = '''4.83'
$ws.Range("E47").Value = '  +8.31%  '

$ws.Range("E49").Value = '  +0.03%  '

$ws.Range("D50").Value = '''1.16'
$ws.Range("E50").Value = '  +6.58%  '

$ws.Range("E51").Value = '  +14.52%  '
